$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 12: a request using an unrecognised term, and the assistant's refusal.
# Write the instruction (column C) before the refusal (column B) so new
# shared-string entries land in the same order as the reference workbook.
$ws.Range("A12").Value = "jsons_train/rapport_original.json"
$ws.Range("C12").Value = "Ajoute un xljfsn du nombre distinct de session_id"
$ws.Range("B12").Value = 'Desole je ne reconnais pas le terme "xljfsn"'

# Row 13: a request for an unsupported Power BI visual, and the refusal
$ws.Range("A13").Value = "jsons_train/rapport_original.json"
$ws.Range("C13").Value = "Ajoute un graphique en etoile des variables de la table"
$ws.Range("B13").Value = 'Desole mais le visuel "graphique en etoile" n''est pas disponible dans Power BI'

# Widen columns A and B, and size the new column C, to fit the longer
# instruction/refusal text now stored in the sheet (matching the saved
# column widths of ~32.73 / ~60.91 / ~82.18 characters, rounded to the
# nearest width this engine's column-width grid can represent).
$ws.Columns.Item(1).ColumnWidth = 31.833333333333332
$ws.Columns.Item(2).ColumnWidth = 60
$ws.Columns.Item(3).ColumnWidth = 81.33333333333333

# Match the selection state recorded in the saved workbook
$ws.Range("B7").Select()
